$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.622.88"
$ws.Range("E2").Value = "  -7.15%  "
$ws.Range("D3").Value = "1.694.61"
$ws.Range("E3").Value = "  -5.74%  "
$ws.Range("E4").Value = "  +0.31%  "
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "219.54"
$c.Style = $s
$ws.Range("E5").Value = "  -5.26%  "
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5119"
$c.Style = $s
$ws.Range("E6").Value = "  -12.84%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -4.50%  "
$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "22.13"
$c.Style = $s
$ws.Range("E9").Value = "  -4.74%  "
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06285"
$c.Style = $s
$ws.Range("E10").Value = "  -7.49%  "
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07364"
$c.Style = $s
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "1.697.03"
$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.518"
$c.Style = $s
$ws.Range("E13").Value = "  -5.75%  "
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5784"
$c.Style = $s
$ws.Range("E14").Value = "  -6.67%  "
$ws.Range("E15").Value = "  -5.65%  "
$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000008423"
$c.Style = $s
$ws.Range("E16").Value = "  -7.73%  "
$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "65.48"
$c.Style = $s
$ws.Range("E17").Value = "  -13.16%  "
$ws.Range("D18").Value = "26.667.96"
$ws.Range("E18").Value = "  -6.93%  "
$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.990"
$c.Style = $s
$ws.Range("E19").Value = "  -8.88%  "
$ws.Range("E20").Value = "  +0.21%  "
$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.98"
$c.Style = $s
$ws.Range("E21").Value = "  -4.69%  "
$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "186.52"
$c.Style = $s
$ws.Range("E22").Value = "  -11.50%  "
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.250"
$c.Style = $s
$ws.Range("E23").Value = "  -8.54%  "
$ws.Range("E24").Value = "  +0.31%  "
$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "144.72"
$c.Style = $s
$ws.Range("E25").Value = "  -5.83%  "
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.497"
$c.Style = $s
$ws.Range("E26").Value = "  -5.87%  "
$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1157"
$c.Style = $s
$ws.Range("E27").Value = "  -8.58%  "
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "15.80"
$c.Style = $s
$ws.Range("E28").Value = "  -3.96%  "
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.358"
$c.Style = $s
$ws.Range("E29").Value = "  -4.42%  "
$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05664"
$c.Style = $s
$ws.Range("E30").Value = "  -7.45%  "
$ws.Range("E31").Value = "  -6.31%  "
$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.507"
$c.Style = $s
$ws.Range("E32").Value = "  -7.45%  "
$ws.Range("E33").Value = "  -8.34%  "
$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.654"
$c.Style = $s
$ws.Range("E34").Value = "  -4.90%  "
$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.019"
$c.Style = $s
$ws.Range("E35").Value = "  -3.00%  "
$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.6005"
$c.Style = $s
$ws.Range("E36").Value = "  -6.76%  "
$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.363"
$c.Style = $s
$ws.Range("E37").Value = "  -5.57%  "
$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.698"
$c.Style = $s
$ws.Range("E38").Value = "  -0.77%  "
$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01616"
$c.Style = $s
$ws.Range("E39").Value = "  -4.92%  "
$ws.Range("D40").Value = "1.102.53"
$ws.Range("E40").Value = "  -3.61%  "
$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8585"
$c.Style = $s
$ws.Range("E41").Value = "  -2.86%  "
$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.834"
$c.Style = $s
$ws.Range("E42").Value = "  -10.34%  "
$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = $s
$ws.Range("E43").Value = "  -0.30%  "
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "99.53"
$c.Style = $s
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Value = "1.853.08"
$ws.Range("E45").Value = "  -5.09%  "
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.00000000114"
$c.Style = $s
$ws.Range("E46").Value = "  +1.62%  "
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "56.51"
$c.Style = $s
$ws.Range("E47").Value = "  -6.20%  "
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = $s
$ws.Range("E48").Value = "  +0.44%  "
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.094"
$c.Style = $s
$ws.Range("E49").Value = "  -2.90%  "
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05241"
$c.Style = $s
$ws.Range("E50").Value = "  -4.45%  "
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4323"
$c.Style = $s
$ws.Range("E51").Value = "  -3.49%  "
